$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update the "Förändrad" date column (C) for rows 2-10 from serial 45224 to 45233
$ws.Range("C2:C10").Value = 45233
